# Updates the crypto price table (columns B-E, rows 2-51) to match the
# latest pull from coinranking.com. Rows 2-8 are simple price/volume
# refreshes; from row 9 on, "OKB" dropped out of the top list causing every
# following coin to shift up one row, with "Aave" newly appearing at the
# bottom (row 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; B = $null; C = $null; D = '28.554.03'; E = '  -1.20%  ' },
    @{ Row = 3; B = $null; C = $null; D = '1.892.50'; E = '  +0.57%  ' },
    @{ Row = 4; B = $null; C = $null; D = '1.009'; E = '  +0.40%  ' },
    @{ Row = 5; B = $null; C = $null; D = '326.66'; E = '  -0.11%  ' },
    @{ Row = 6; B = $null; C = $null; D = '1.008'; E = $null },
    @{ Row = 7; B = $null; C = $null; D = '0.4592'; E = '  -1.55%  ' },
    @{ Row = 8; B = $null; C = $null; D = '0.3864'; E = '  -2.32%  ' },
    @{ Row = 9; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.07870'; E = '  -0.76%  ' },
    @{ Row = 10; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '1.002'; E = '  +2.48%  ' },
    @{ Row = 11; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '21.64'; E = '  -3.54%  ' },
    @{ Row = 12; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.909.07'; E = '  +3.78%  ' },
    @{ Row = 13; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '7.085'; E = '  +1.21%  ' },
    @{ Row = 14; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.716'; E = '  -0.88%  ' },
    @{ Row = 15; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.06965'; E = '  -0.36%  ' },
    @{ Row = 16; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '87.47'; E = '  -1.62%  ' },
    @{ Row = 17; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.009'; E = '  +0.48%  ' },
    @{ Row = 18; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.00001005'; E = '  -1.04%  ' },
    @{ Row = 19; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '17.22'; E = '  +1.07%  ' },
    @{ Row = 20; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.007'; E = '  +0.35%  ' },
    @{ Row = 21; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '28.582.26'; E = '  -1.09%  ' },
    @{ Row = 22; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '5.326'; E = '  -0.68%  ' },
    @{ Row = 23; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '11.01'; E = '  -1.10%  ' },
    @{ Row = 24; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.127.51'; E = '  +2.61%  ' },
    @{ Row = 25; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '2.061'; E = '  -2.61%  ' },
    @{ Row = 26; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '154.84'; E = '  +0.91%  ' },
    @{ Row = 27; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '19.40'; E = '  -0.35%  ' },
    @{ Row = 28; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '5.855'; E = '  +1.34%  ' },
    @{ Row = 29; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '1.960'; E = '  -2.64%  ' },
    @{ Row = 30; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '118.32'; E = '  -1.39%  ' },
    @{ Row = 31; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.09327'; E = '  -0.77%  ' },
    @{ Row = 32; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.9251'; E = '  -2.24%  ' },
    @{ Row = 33; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '5.298'; E = '  -0.69%  ' },
    @{ Row = 34; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.336'; E = '  -1.40%  ' },
    @{ Row = 35; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '3.270'; E = '  -2.36%  ' },
    @{ Row = 36; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.05766'; E = '  -2.76%  ' },
    @{ Row = 37; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '1.155'; E = '  +0.19%  ' },
    @{ Row = 38; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.02073'; E = '  -2.37%  ' },
    @{ Row = 39; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '7.795'; E = '  -1.92%  ' },
    @{ Row = 40; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.5671'; E = '  -1.10%  ' },
    @{ Row = 41; B = 'Algorand'; C = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D = '0.1789'; E = '  -0.48%  ' },
    @{ Row = 42; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '9.752'; E = '  -2.60%  ' },
    @{ Row = 43; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '11.82'; E = '  -0.26%  ' },
    @{ Row = 44; B = 'Cronos'; C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = '0.07166'; E = '  -1.10%  ' },
    @{ Row = 45; B = 'Decentraland'; C = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D = '0.5356'; E = '  +0.01%  ' },
    @{ Row = 46; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '2.200'; E = '  +2.76%  ' },
    @{ Row = 47; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '1.842'; E = '  -0.74%  ' },
    @{ Row = 48; B = 'WEMIXToken'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '1.116'; E = '  -2.01%  ' },
    @{ Row = 49; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '112.78'; E = '  -1.24%  ' },
    @{ Row = 50; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.466'; E = '  +4.28%  ' },
    @{ Row = 51; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '71.00'; E = '  +1.74%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.B) { $ws.Range("B$row").Value = $u.B }
    if ($null -ne $u.C) { $ws.Range("C$row").Value = $u.C }

    if ($null -ne $u.D) {
        $cell = $ws.Range("D$row")
        # Price strings use a dotted/grouped format (e.g. "28.554.03",
        # "1.009", "0.00001005") that Excel would otherwise silently
        # reinterpret as a number (dropping meaningful trailing/leading
        # zeros or flipping to scientific notation). Force a text entry
        # with a leading apostrophe, then strip the resulting "Text"
        # number-format back off so the cell keeps its original (default)
        # style.
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }

    if ($null -ne $u.E) { $ws.Range("E$row").Value = $u.E }
}
